# =========================================================================
# Edit: add a new "storage_content__Building1" sheet (with per-storage
# content series), drop the now-redundant storage_content column from
# shSourceBus__Building1, fix the row order of storage related labels on
# env_impacts__Building1 / capStorages__Building1, and refresh a couple of
# recomputed electricityBus__Building1 values.
# =========================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1. Create the new "storage_content__Building1" sheet, positioned right
#    before "costs__Building1" (i.e. as the 9th sheet).
#    We clone shSourceBus__Building1 because it already has the correct
#    date column (A) + header/date cell styles, then strip it down to
#    just the two storage-content columns.
# -------------------------------------------------------------------
$costsSheet = $wb.Worksheets.Item("costs__Building1")
$srcBusSheet = $wb.Worksheets.Item("shSourceBus__Building1")
$srcBusSheet.Copy($costsSheet)
$newSheet = $wb.Worksheets.Item("shSourceBus__Building1 (2)")
$newSheet.Name = "storage_content__Building1"

# Original columns on the clone: B=HP->shSourceBus flow, C=shSourceBus->shSource
# flow, D=shSourceBus->shStorage flow, E=storage_content (shStorage content).
# Remove the three flow columns (B, C, D) so old column E becomes column B,
# keeping its header/values (shStorage storage content) intact.
$newSheet.Columns.Item(4).Delete()
$newSheet.Columns.Item(3).Delete()
$newSheet.Columns.Item(2).Delete()

# Rename the remaining header and add the electricalStorage content column.
$newSheet.Range("B1").Value = "shStorage__B001_storage_content"
$newSheet.Range("C1").Value = "electricalStorage__B001_storage_content"

# C1 needs the same bold/border header styling as B1 (copy format only).
$newSheet.Range("B1").Copy()
$newSheet.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# electricalStorage__Building1 storage-content values (column C), rows 2-26.
$elecStorageContent = @(
    119.4539202662791,
    112.081874825057,
    111.4428216692431,
    110.7141007390105,
    106.3290519046023,
    99.11654483246862,
    92.68741784255559,
    85.81855423530124,
    56.98778943746109,
    56.98778943746109,
    56.98778943746109,
    56.98778943746109,
    56.98778943746109,
    53.70285865349751,
    49.96739538562801,
    48.34128242865126,
    46.20251166469777,
    43.87013624958149,
    33.81247844260475,
    19.72698010539545,
    4.221096345348838,
    2.256411960465116,
    1.310282391860465,
    0.6390863790697674,
    0
)
for ($i = 0; $i -lt $elecStorageContent.Count; $i++) {
    $row = $i + 2
    $newSheet.Cells.Item($row, 3).Value = $elecStorageContent[$i]
}

# -------------------------------------------------------------------
# 2. shSourceBus__Building1: the storage_content column (E) has moved to
#    its own sheet, so drop it here.
# -------------------------------------------------------------------
$srcBusSheet.Columns.Item(5).Delete()

# -------------------------------------------------------------------
# 3. electricityBus__Building1: refreshed MPC re-solve values on rows 10
#    and 15 (columns B and C only, D is unchanged).
# -------------------------------------------------------------------
$elecBusSheet = $wb.Worksheets.Item("electricityBus__Building1")
$elecBusSheet.Range("B10").Value = 24.79445772614253
$elecBusSheet.Range("C10").Value = 21.6502397029283
$elecBusSheet.Range("B15").Value = 2.825040474208677
$elecBusSheet.Range("C15").Value = 0

# -------------------------------------------------------------------
# 4. env_impacts__Building1: swap the electricalStorage / shStorage rows
#    (rows 6 and 7) so shStorage comes first.
# -------------------------------------------------------------------
$envSheet = $wb.Worksheets.Item("env_impacts__Building1")
$a6 = $envSheet.Range("A6").Value()
$b6 = $envSheet.Range("B6").Value()
$a7 = $envSheet.Range("A7").Value()
$b7 = $envSheet.Range("B7").Value()
$envSheet.Range("A6").Value = $a7
$envSheet.Range("B6").Value = $b7
$envSheet.Range("A7").Value = $a6
$envSheet.Range("B7").Value = $b6

# -------------------------------------------------------------------
# 5. capStorages__Building1: swap the electricalStorage / shStorage rows
#    (rows 2 and 3) so shStorage comes first.
# -------------------------------------------------------------------
$capStorSheet = $wb.Worksheets.Item("capStorages__Building1")
$a2 = $capStorSheet.Range("A2").Value()
$b2 = $capStorSheet.Range("B2").Value()
$a3 = $capStorSheet.Range("A3").Value()
$b3 = $capStorSheet.Range("B3").Value()
$capStorSheet.Range("A2").Value = $a3
$capStorSheet.Range("B2").Value = $b3
$capStorSheet.Range("A3").Value = $a2
$capStorSheet.Range("B3").Value = $b2

# -------------------------------------------------------------------
# 6. Restore the original active tab (gridBus__Building1) so the copy
#    operations above don't leave a different sheet selected.
# -------------------------------------------------------------------
$wb.Worksheets.Item(1).Select()
